$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 362
$ws1.Range("F3").Value = 66
$ws1.Range("F4").Value = 274
$ws1.Range("F5").Value = 4088

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 362
$ws4.Range("F3").Value = 66
$ws4.Range("F4").Value = 274
$ws4.Range("F5").Value = 4088
